$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices / 1h volume deltas), matching the source
# feed pull. A couple of coins also swapped rank position (rows 26/27 and
# 36/37), so their Coin name + Link are rewritten too.
#
# Numeric-looking price strings (e.g. "1.00", "507.91") are written with a
# leading apostrophe so Excel keeps them as text (matching the workbook's
# existing inline-string price format) instead of coercing them to numbers;
# the style is then reset to Normal so no stray numeric format sticks to
# the cell.

$ws.Range("D2").Value = '56.376.10'
$ws.Range("E2").Value = '  -1.49%  '

$ws.Range("D3").Value = '3.008.74'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = "'507.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").Value = "'138.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.18%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = "'7.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.35%  '

$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("E11").Value = '  +3.30%  '

$ws.Range("D12").Value = '3.518.42'
$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").Value = "'25.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").Value = "'0.0000162"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.86%  '

$ws.Range("D16").Value = '56.328.51'

$ws.Range("D17").Value = '3.005.48'
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").Value = "'5.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.88%  '

$ws.Range("D19").Value = "'12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.76%  '

$ws.Range("D20").Value = "'8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.01%  '

$ws.Range("D21").Value = "'333.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.24%  '

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = "'0.498"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("D24").Value = "'64.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.71%  '

$ws.Range("D25").Value = '3.128.65'
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.53%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '

$ws.Range("D28").Value = '0.0₃0939'
$ws.Range("E28").Value = '  +5.46%  '

$ws.Range("E29").Value = '  -3.82%  '

$ws.Range("E30").Value = '  -3.25%  '

$ws.Range("D31").Value = "'1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '

$ws.Range("D32").Value = "'20.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").Value = "'152.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.71%  '

$ws.Range("D35").Value = "'4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.20%  '

$ws.Range("B36").Value = 'EnergySwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D36").Value = "'26.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.19%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'5.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("E38").Value = '  -1.01%  '

$ws.Range("D39").Value = "'0.0662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '

$ws.Range("D40").Value = '3.047.03'
$ws.Range("E40").Value = '  +0.81%  '

$ws.Range("D41").Value = "'36.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.80%  '

$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("D43").Value = "'3.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.45%  '

$ws.Range("D44").Value = "'0.655"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.45%  '

$ws.Range("D45").Value = '2.200.94'
$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("E46").Value = '  -2.56%  '

$ws.Range("E47").Value = '  +2.27%  '

$ws.Range("D48").Value = "'0.925"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '

$ws.Range("D49").Value = "'5.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.73%  '

$ws.Range("D50").Value = "'19.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.19%  '

$ws.Range("D51").Value = "'0.0851"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.95%  '
